$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 138 (shifts existing rows 138-154 down to 139-155)
$ws.Rows("138:138").Insert()

# Populate the newly inserted row 138 with the new weekly record
$ws.Range("A138").Value = 7
$ws.Range("B138").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C138").Value = "Ñuble"
$ws.Range("D138").Value = 44449
$ws.Range("E138").Value = 16
$ws.Range("F138").Value = 100112023
$ws.Range("G138").Value = "Brócoli"
$ws.Range("H138").Value = "Sin especificar"
$ws.Range("I138").Value = "Primera"
$ws.Range("J138").Value = 160
$ws.Range("K138").Value = 700
$ws.Range("L138").Value = 750
$ws.Range("M138").Value = 725
$ws.Range("N138").Value = "$/unidad"
$ws.Range("O138").Value = "Región del Maule"
$ws.Range("P138").Value = 725
$ws.Range("Q138").Value = 1
$ws.Range("R138").Value = "Hortaliza"
